# edit.ps1 - apply "carotenuto first draft complete" revision
#
# Summary of changes (see unified diff):
#   1. Title paragraph gains a leading ", " run before
#      "Effectiveness of Rattlesnake Aversion Training in Dogs".
#   2. The "Your animal has undergone..." (Survey 1) paragraph: several
#      runs that used to be split get merged into one run, and the
#      " helpful for future dog owners..." run gets split so that the
#      word "helpful" is wrapped in gramStart/gramEnd proofErr marks.
#   3. "Age, breed, sex of dog now" paragraph: the trailing " " run and
#      "(Drop down menu)" run merge into a single " (Drop down menu)" run.
#   4. The "Your animal was recently envenomated..." (Survey 2) paragraph:
#      the " helpful for future dog owners" run gets split the same way
#      as in change 2 (gramStart/gramEnd around "helpful"), while the
#      rest of the paragraph's runs stay as they were.
#
# Because Word normalizes (merges) adjacent runs that end up with
# identical formatting when a document is saved, the most reliable way
# to reproduce an exact target run layout is to rebuild each affected
# paragraph's full contents in one shot via Range.InsertXML, rather than
# trying to surgically edit individual runs in place.

$d = $word.ActiveDocument

function Get-ParagraphByText($doc, [string]$containsText) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs($i)
        if ($p.Range.Text.Contains($containsText)) {
            return $p
        }
    }
    throw "Paragraph containing '$containsText' not found"
}

function Set-ParagraphXml($paragraph, [string]$innerBodyXml) {
    $rng = $paragraph.Range
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' + $innerBodyXml + '</w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($pkg)
}

# ---------------------------------------------------------------------
# 1. Title: insert ", " as its own run before the title text.
# ---------------------------------------------------------------------
$titlePara = Get-ParagraphByText $d "Effectiveness of Rattlesnake Aversion Training in Dogs"
$titleXml = '<w:p>' +
              '<w:r><w:t xml:space="preserve">, </w:t></w:r>' +
              '<w:r><w:t>Effectiveness of Rattlesnake Aversion Training in Dogs</w:t></w:r>' +
            '</w:p>'
Set-ParagraphXml $titlePara $titleXml

# ---------------------------------------------------------------------
# 2. Survey 1 paragraph: "Your animal has undergone ..."
# ---------------------------------------------------------------------
$survey1Para = Get-ParagraphByText $d "Your animal has undergone"
$survey1Xml = '<w:p>' +
    '<w:r><w:t xml:space="preserve">Your animal has undergone rattlesnake avoidance training with Rattlesnake Ready. The University of Arizona College of Veterinary Medicine is conducting a survey of owners. If you would kindly take 15 minutes to participate in this survey, it </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>would</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>helpful</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> for future dog owners. Two owners will be drawn at random for a $50 amazon gift card. </w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml $survey1Para $survey1Xml

# ---------------------------------------------------------------------
# 3. "Age, breed, sex of dog now (Drop down menu)" paragraph: merge the
#    trailing " " + "(Drop down menu)" runs into a single run.
# ---------------------------------------------------------------------
$ageNowPara = Get-ParagraphByText $d "Age, breed, sex of dog now"
$ageNowXml = '<w:p>' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Age, breed, sex of dog now</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> (Drop down menu)</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml $ageNowPara $ageNowXml

# ---------------------------------------------------------------------
# 4. Survey 2 paragraph: "Your animal was recently envenomated ..."
#    Only the " helpful for future dog owners" run is split (gramStart/
#    gramEnd around "helpful"); the rest of the runs are left untouched.
# ---------------------------------------------------------------------
$survey2Para = Get-ParagraphByText $d "Your animal was recently envenomated"
$survey2Xml = '<w:p>' +
    '<w:r><w:t>Your animal was recently envenomated by a rattlesnake and you sought care at a local emergency practice. The University of Arizona</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> College of Veterinary Medicine</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> is conduct</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">ing a survey of owners. If you </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">would kindly </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">take 5 minutes to participate in this survey, it </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>would</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>helpful</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> for future dog owners</w:t></w:r>' +
    '<w:r><w:t>. Two</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> owners will be drawn </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">at random </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">for a $50 amazon gift card. </w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml $survey2Para $survey2Xml

Write-Host "All edits applied."
